$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Update column B (2024 performance) values per row, per the diff
$ws.Cells.Item(2, 2).Value = -0.11145510835913321
$ws.Cells.Item(3, 2).Value = 0.1080213849698837
$ws.Cells.Item(4, 2).Value = 0.1064854759850451
$ws.Cells.Item(5, 2).Value = -0.054995970991135963
$ws.Cells.Item(6, 2).Value = 0.065626681011295762
$ws.Cells.Item(7, 2).Value = 0.098613251155623027
$ws.Cells.Item(8, 2).Value = 0.16674197921373629
$ws.Cells.Item(9, 2).Value = 0.06572144126676438
$ws.Cells.Item(10, 2).Value = 0.1024925198965982
$ws.Cells.Item(11, 2).Value = 0.15737473535638569
$ws.Cells.Item(12, 2).Value = 0.20844055584148211
$ws.Cells.Item(13, 2).Value = -0.1035707779194754
$ws.Cells.Item(14, 2).Value = 0.1512388162422578
$ws.Cells.Item(15, 2).Value = -0.01020408163265518
$ws.Cells.Item(16, 2).Value = -0.099429115128448764
$ws.Cells.Item(17, 2).Value = 0.33494872563012312
$ws.Cells.Item(18, 2).Value = -0.04953497775980531
$ws.Cells.Item(19, 2).Value = 0.062525375558264429
$ws.Cells.Item(20, 2).Value = 0.13551401869158911
$ws.Cells.Item(21, 2).Value = 0.053673627223510767
$ws.Cells.Item(22, 2).Value = -0.077795104835467588
$ws.Cells.Item(23, 2).Value = 0.15661252900232039
$ws.Cells.Item(24, 2).Value = 0.1960000000000004
$ws.Cells.Item(25, 2).Value = 0.1978319783197833
$ws.Cells.Item(26, 2).Value = 0.15695346795434631
$ws.Cells.Item(27, 2).Value = 0.02507712210170188
$ws.Cells.Item(28, 2).Value = 0.15339902204943151
$ws.Cells.Item(29, 2).Value = 0.2442384769539061
$ws.Cells.Item(30, 2).Value = 0.23721954948069121
$ws.Cells.Item(31, 2).Value = 0.081820050709243919
$ws.Cells.Item(32, 2).Value = 0.1209915539380626
$ws.Cells.Item(33, 2).Value = 0.058163720101025218
$ws.Cells.Item(34, 2).Value = 0.088006986221619732
$ws.Cells.Item(35, 2).Value = -0.077976817702844148
$ws.Cells.Item(36, 2).Value = 0.11336982017200881
$ws.Cells.Item(37, 2).Value = -0.054995970991135963
$ws.Cells.Item(38, 2).Value = -0.1017942145734156
$ws.Cells.Item(39, 2).Value = 0.1949567181031244
$ws.Cells.Item(40, 2).Value = -0.08941485864562837
$ws.Cells.Item(41, 2).Value = 0.042255511588468853
$ws.Cells.Item(42, 2).Value = 0.2035928143712564
$ws.Cells.Item(43, 2).Value = 0.12836624775583449
$ws.Cells.Item(44, 2).Value = 0.2345554195711155
$ws.Cells.Item(45, 2).Value = 0.072084160807257769
$ws.Cells.Item(46, 2).Value = 0.1041515517936311
$ws.Cells.Item(47, 2).Value = 0.23921683734878291
$ws.Cells.Item(48, 2).Value = -0.11759504862953141
$ws.Cells.Item(49, 2).Value = 0.23151645979492949
$ws.Cells.Item(50, 2).Value = 0.062525375558264429
$ws.Cells.Item(51, 2).Value = 0.13508260447036119
$ws.Cells.Item(52, 2).Value = 0.044333149601808897
$ws.Cells.Item(53, 2).Value = 0.25660226561956878
$ws.Cells.Item(54, 2).Value = 0.2896855398598932
$ws.Cells.Item(55, 2).Value = 0.2605398675796502
$ws.Cells.Item(56, 2).Value = 0.1613361762615482
$ws.Cells.Item(57, 2).Value = 0.11354817140878801
$ws.Cells.Item(58, 2).Value = 0.1122944452457038
$ws.Cells.Item(60, 2).Value = 0.28316197539187637
$ws.Cells.Item(61, 2).Value = 0.10280569514237679
$ws.Cells.Item(62, 2).Value = 0.2099832211711701
$ws.Cells.Item(63, 2).Value = 0.26649041375039179
$ws.Cells.Item(64, 2).Value = 0.43636363636363562
$ws.Cells.Item(65, 2).Value = 0.19541875447387369
$ws.Cells.Item(66, 2).Value = 0.33487677537260557
$ws.Cells.Item(67, 2).Value = 0.2512421815631023
$ws.Cells.Item(68, 2).Value = 0.2209185561630633
$ws.Cells.Item(69, 2).Value = 0.27036245229021172
$ws.Cells.Item(70, 2).Value = 0.24389314105734991
$ws.Cells.Item(71, 2).Value = 0.21189206832771279
$ws.Cells.Item(72, 2).Value = 0.22918654464200361
$ws.Cells.Item(73, 2).Value = 0.23914592902533191
$ws.Cells.Item(74, 2).Value = 0.14163017671182329
$ws.Cells.Item(75, 2).Value = 0.1236887844102068
$ws.Cells.Item(76, 2).Value = 0.1598281700848361
$ws.Cells.Item(77, 2).Value = 0.18489055269588769
$ws.Cells.Item(78, 2).Value = 0.20438930024681559
$ws.Cells.Item(79, 2).Value = 0.1211136973086855

# Update the active cell selection shown on the sheet view
$ws.Range("I73").Select()
